$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H29").Value = 1417.9
$ws.Range("I29").Value = 135.8
$ws.Range("J29").Value = 2700
$ws.Range("K29").Value = 407.4
$ws.Range("L29").Value = 8100
$ws.Range("M29").Value = -126.4
$ws.Range("N29").Value = -8662

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H38").Value = 1321.762
$ws.Range("I38").Value = 63.083332
$ws.Range("K38").Value = 189.249996
$ws.Range("M38").Value = 182.750004

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H58").Value = 1203.4375
$ws.Range("I58").Value = 479.58334
$ws.Range("J58").Value = 3375
$ws.Range("K58").Value = 1438.75002
$ws.Range("L58").Value = 10125
$ws.Range("M58").Value = -1288.75002
$ws.Range("N58").Value = -10425

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H80").Value = 2485.1428
$ws.Range("I80").Value = 1318.4
$ws.Range("J80").Value = 3133.3333
$ws.Range("K80").Value = 3955.2
$ws.Range("L80").Value = 9399.999899999999
$ws.Range("M80").Value = -2957.2
$ws.Range("N80").Value = -11395.9999

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H83").Value = 2485.1428
$ws.Range("I83").Value = 1318.4
$ws.Range("J83").Value = 3133.3333
$ws.Range("K83").Value = 11865.6
$ws.Range("L83").Value = 28199.9997
$ws.Range("M83").Value = -6873.6
$ws.Range("N83").Value = -38183.9997

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H125").Value = 1580.6666
$ws.Range("I125").Value = 1216
$ws.Range("J125").Value = 1653.6
$ws.Range("K125").Value = 10944
$ws.Range("L125").Value = 14882.4
$ws.Range("M125").Value = -8484
$ws.Range("N125").Value = -19802.4

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H127").Value = 811.38464
$ws.Range("J127").Value = 1430.25
$ws.Range("L127").Value = 4290.75
$ws.Range("N127").Value = -14210.75

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H133").Value = 58316.668
$ws.Range("J133").Value = 58316.668
$ws.Range("L133").Value = 58316.668
$ws.Range("N133").Value = -68436.66800000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5417.72
$ws.Range("I32").Value = 4244.2676
$ws.Range("J32").Value = 12626.071
$ws.Range("K32").Value = 4244.2676
$ws.Range("L32").Value = 12626.071
$ws.Range("M32").Value = -3957.2676
$ws.Range("N32").Value = -13200.071

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 279131.8
$ws.Range("I134").Value = 385266.7
$ws.Range("J134").Value = 3181.1
$ws.Range("K134").Value = 1155800.1
$ws.Range("L134").Value = 9543.299999999999
$ws.Range("M134").Value = -1153265.1
$ws.Range("N134").Value = -14613.3

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 2191.4211
$ws.Range("I132").Value = 1821.1875
$ws.Range("J132").Value = 4166
$ws.Range("K132").Value = 5463.5625
$ws.Range("L132").Value = 12498
$ws.Range("M132").Value = -2933.5625
$ws.Range("N132").Value = -17558

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H6").Value = 66667108
$ws.Range("I6").Value = 83333500
$ws.Range("K6").Value = 250000500
$ws.Range("M6").Value = -250000387

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 1362.1794
$ws.Range("I113").Value = 465.73077
$ws.Range("J113").Value = 3155.077
$ws.Range("K113").Value = 1397.19231
$ws.Range("L113").Value = 9465.231
$ws.Range("M113").Value = 772.8076900000001
$ws.Range("N113").Value = -13805.231

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H141").Value = 1343.9333
$ws.Range("I141").Value = 1031.8
$ws.Range("J141").Value = 1500
$ws.Range("K141").Value = 3095.4
$ws.Range("L141").Value = 4500
$ws.Range("M141").Value = 2084.6
$ws.Range("N141").Value = -14860

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H43").Value = 3489.1155
$ws.Range("I43").Value = 1165.1818
$ws.Range("J43").Value = 5193.3335
$ws.Range("K43").Value = 1165.1818
$ws.Range("L43").Value = 5193.3335
$ws.Range("M43").Value = -1014.1818
$ws.Range("N43").Value = -5495.3335

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H57").Value = 14430.637
$ws.Range("J57").Value = 15073.7
$ws.Range("L57").Value = 15073.7
$ws.Range("N57").Value = -16713.7

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1125.5483
$ws.Range("I102").Value = 1000.5217
$ws.Range("J102").Value = 1485
$ws.Range("K102").Value = 1000.5217
$ws.Range("L102").Value = 1485
$ws.Range("M102").Value = 621.4783
$ws.Range("N102").Value = -4729

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 2884.9473
$ws.Range("I126").Value = 1601.0769
$ws.Range("J126").Value = 5666.6665
$ws.Range("K126").Value = 4803.2307
$ws.Range("L126").Value = 16999.9995
$ws.Range("M126").Value = -2333.2307
$ws.Range("N126").Value = -21939.9995

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2564878.2
$ws.Range("I22").Value = 8333503
$ws.Range("J22").Value = 1045.1111
$ws.Range("K22").Value = 8333503
$ws.Range("L22").Value = 1045.1111
$ws.Range("M22").Value = -8333208
$ws.Range("N22").Value = -1635.1111

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H27").Value = 2564878.2
$ws.Range("I27").Value = 8333503
$ws.Range("J27").Value = 1045.1111
$ws.Range("K27").Value = 8333503
$ws.Range("L27").Value = 1045.1111
$ws.Range("M27").Value = -8333396
$ws.Range("N27").Value = -1259.1111

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 281.69696
$ws.Range("I55").Value = 350.75
$ws.Range("J55").Value = 216.70589
$ws.Range("K55").Value = 350.75
$ws.Range("L55").Value = 216.70589
$ws.Range("M55").Value = -177.75
$ws.Range("N55").Value = -562.70589

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 4057.3215
$ws.Range("I122").Value = 4190.7617
$ws.Range("J122").Value = 3657
$ws.Range("K122").Value = 12572.2851
$ws.Range("L122").Value = 10971
$ws.Range("M122").Value = -10122.2851
$ws.Range("N122").Value = -15871

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 12828580
$ws.Range("I132").Value = 20844168
$ws.Range("J132").Value = 3641.2
$ws.Range("K132").Value = 62532504
$ws.Range("L132").Value = 10923.6
$ws.Range("M132").Value = -62529974
$ws.Range("N132").Value = -15983.6

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 14052.929
$ws.Range("J54").Value = 14052.929
$ws.Range("L54").Value = 14052.929
$ws.Range("N54").Value = -15092.929

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H94").Value = 22337.5
$ws.Range("J94").Value = 22337.5
$ws.Range("L94").Value = 22337.5
$ws.Range("N94").Value = -24139.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1440.1333
$ws.Range("I132").Value = 773.2222
$ws.Range("J132").Value = 2440.5
$ws.Range("K132").Value = 2319.6666
$ws.Range("L132").Value = 7321.5
$ws.Range("M132").Value = 210.3334
$ws.Range("N132").Value = -12381.5

Write-Host "Applied all updates"